# Rename header columns B1/C1 from "Prin ID"/"Prin Description"
# to "Partner ID"/"Partner Description".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Partner ID"
$ws.Range("C1").Value = "Partner Description"

# Move the active selection to J11 (single cell), matching the saved view state.
$ws.Range("J11").Select()
